function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook

$wsAdmin = $wb.Worksheets.Item(1)    # Administrativos
$wsDoc   = $wb.Worksheets.Item(2)    # Docentes-Conciliadores
$wsEst   = $wb.Worksheets.Item(3)    # Estudiantes

$lightBlue = RGBColor 183 222 232    # FFB7DEE8 - conditional format highlight
$white     = RGBColor 255 255 255

# ---------------------------------------------------------------------------
# Estudiantes (sheet3): remove the Goku/Saiyajin row, update remaining row,
# change tab color, add a hyperlink + conditional formatting.
# ---------------------------------------------------------------------------
$wsEst.Rows.Item(3).Delete() | Out-Null

$wsEst.Range("A2").Value = "Krilin"
$wsEst.Range("B2").Value = "Aprendiz"
$wsEst.Range("C2").Value = 1298323
$wsEst.Range("D2").Value = 8329842
$wsEst.Range("E2").Value = "krilin@hotmail.com"
$wsEst.Hyperlinks.Add($wsEst.Range("E2"), "mailto:krilin@hotmail.com") | Out-Null

$wsEst.Tab.Color = RGBColor 92 158 49   # FF5C9E31

$hdrEst = $wsEst.Range("A1:E1")
$hdrEst.Font.Color = $white
$hdrEst.Interior.Color = RGBColor 92 158 49

$fcEst = $wsEst.Range("A2:E41").FormatConditions.Add(2, 0, "LEN(TRIM(A2))=0")
$fcEst.Interior.Color = $lightBlue

# ---------------------------------------------------------------------------
# Docentes-Conciliadores (sheet2): update row, drop the Tarjeta_Profesinal
# value, change tab color, update hyperlink, add conditional formatting.
# ---------------------------------------------------------------------------
$wsDoc.Range("A2").Value = "Roshi"
$wsDoc.Range("B2").Value = "Sayajin"
$wsDoc.Range("C2").Value = 123456789
$wsDoc.Range("D2").Value = 12345675
$wsDoc.Range("E2").Value = "goten@ugc.edu.co"
$wsDoc.Range("F2").Clear() | Out-Null

$wsDoc.Range("E2").Hyperlinks.Delete() | Out-Null
$wsDoc.Hyperlinks.Add($wsDoc.Range("E2"), "mailto:goten@ugc.edu.co") | Out-Null

$wsDoc.Tab.Color = RGBColor 0 70 15     # FF00460F

$hdrDoc = $wsDoc.Range("A1:F1")
$hdrDoc.Font.Color = $white
$hdrDoc.Interior.Color = RGBColor 0 70 15

$fcDoc = $wsDoc.Range("A2:E41").FormatConditions.Add(2, 0, "LEN(TRIM(A2))=0")
$fcDoc.Interior.Color = $lightBlue

# ---------------------------------------------------------------------------
# Administrativos (sheet1): update row, update hyperlink, add conditional
# formatting. Tab color / header color stay the same green.
# ---------------------------------------------------------------------------
$wsAdmin.Range("A2").Value = "Son "
$wsAdmin.Range("B2").Value = "Goku"
$wsAdmin.Range("C2").Value = 1249239
$wsAdmin.Range("D2").Value = 314439376
$wsAdmin.Range("E2").Value = "son.goku@ugc.edu.co"

$wsAdmin.Range("E2").Hyperlinks.Delete() | Out-Null
$wsAdmin.Hyperlinks.Add($wsAdmin.Range("E2"), "mailto:son.goku@ugc.edu.co") | Out-Null

$hdrAdmin = $wsAdmin.Range("A1:E1")
$hdrAdmin.Font.Color = $white
$hdrAdmin.Interior.Color = RGBColor 0 145 61   # FF00913D

$fcAdmin = $wsAdmin.Range("A2:E41").FormatConditions.Add(2, 0, "LEN(TRIM(A2))=0")
$fcAdmin.Interior.Color = $lightBlue

# ---------------------------------------------------------------------------
# Selections per sheet (select in tab order, Estudiantes last so it stays the
# active tab), matching the final workbook selection state.
# ---------------------------------------------------------------------------
$wsAdmin.Range("B3").Select() | Out-Null
$wsDoc.Range("A2").Select() | Out-Null
$wsEst.Range("D3").Select() | Out-Null
